$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 471.4767337671348
$ws.Range("B2").Value = 4679.739
$ws.Range("C2").Value = -1121.524

# Delete rows 3 and 4 (entire rows) so the used range shrinks to A1:C2
$ws.Range("A3:C4").EntireRow.Delete()
